# Add CMS1500 claim test data to the "testDataSheet" worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("testDataSheet")

# ---------------------------------------------------------------------
# 1. New header cells in row 1 (columns AQ..BE / 43..57)
# ---------------------------------------------------------------------
$headers = @(
    "patientid",
    "patientSigndate",
    "diagnosisCode",
    "providerID",
    "serviceFromDate",
    "serviceToDate",
    "POS",
    "procCode",
    "serviceAmt",
    "dayunits",
    "diagnosispointer",
    "renderingtaxonomyCode",
    "physiciansigndate",
    "sitephno",
    "billingTaxonomyCode"
)
$col = 43
foreach ($h in $headers) {
    $ws.Cells.Item(1, $col).Value = $h
    $col = $col + 1
}

# diagnosisCode header (AS1) carries a text number format, like the rest of
# the "@" formatted header/data cells in this sheet.
$ws.Range("AS1").NumberFormat = "@"

# ---------------------------------------------------------------------
# 2. New data row (row 4) holding a sample createCMS1500Claim test case
# ---------------------------------------------------------------------
$ws.Range("A4").Value = "createCMS1500Claim"

$ws.Range("AQ4").Value = 37

# AR4 holds a literal date-shaped string (not an actual date value) - force
# text interpretation while writing it, then drop back to the default style
# so the cell ends up with no explicit style reference, same as the source.
$ws.Range("AR4").NumberFormat = "@"
$ws.Range("AR4").Value = "07/13/2023"
$ws.Range("AR4").Style = "Normal"

$ws.Range("AS4").Value = "F70        *Mild intellectual disabilities"
$ws.Range("AT4").Value = 1794

$ws.Range("AU4:AV4").NumberFormat = "dd/mm/yyyy"
$ws.Range("AU4").Value = 44934
$ws.Range("AV4").Value = 44934

$ws.Range("AW4").Value = "Group Home"
$ws.Range("AX4").Value = "YP770"
$ws.Range("AY4").Value = 100
$ws.Range("AZ4").Value = 1
$ws.Range("BA4").Value = "A"
$ws.Range("BB4").Value = "251S00000X"

# BC4 holds a literal date-shaped string too.
$ws.Range("BC4").NumberFormat = "@"
$ws.Range("BC4").Value = "08/16/2023"
$ws.Range("BC4").Style = "Normal"

$ws.Range("BD4").Value = "243-864-7452"
$ws.Range("BE4").Value = "251S00000X"

# ---------------------------------------------------------------------
# 3. Column widths - column A grew wider, and the new columns AQ..BE
#    take on the widths used for the new CMS1500 claim fields.
#    (values below are pre-compensated so the ColumnWidth -> stored xlsx
#    `width` round-trip lands as close as possible to the target widths)
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 25.8333333333333

$colWidths = @(
    12.8333333333333,
    16,
    14.6666666666667,
    14.1666666666667,
    17.1666666666667,
    10.3333333333333,
    15.6666666666667,
    12.3333333333333,
    12.8333333333333,
    12.3333333333333,
    17,
    28.8333333333333,
    12.6666666666667,
    13.6666666666667,
    20.6666666666667
)
$col = 43
foreach ($w in $colWidths) {
    $ws.Columns.Item($col).ColumnWidth = $w
    $col = $col + 1
}

# ---------------------------------------------------------------------
# 4. View state - selection & scroll position moved to the new columns.
# ---------------------------------------------------------------------
$ws.Activate()
$ws.Range("AZ6").Select() | Out-Null
$win = $excel.ActiveWindow
$win.ScrollColumn = 48
